# Updates cryptocurrency price/volume data to the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.800.15'
$ws.Range('E2').Value = '  -0.86%  '
$ws.Range('D3').Value = '1.625.43'
$ws.Range('E3').Value = '  -0.94%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.02'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.33%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5109'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.99%  '
$ws.Range('E7').Value = '  +0.16%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2579'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.35%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06374'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.12%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.35'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.71%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07770'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.07%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.253'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.01%  '
$ws.Range('D13').Value = '1.625.07'
$ws.Range('E13').Value = '  -0.83%  '
$ws.Range('D14').Value = '1.848.69'
$ws.Range('E14').Value = '  -1.03%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5550'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.28%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '63.57'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.63%  '
$ws.Range('D17').Value = '0.0₅7526'
$ws.Range('E17').Value = '  -3.48%  '
$ws.Range('D18').Value = '25.773.57'
$ws.Range('E18').Value = '  -0.98%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.003'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.17%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '193.71'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.86%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.329'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -3.09%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.783'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.16%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.989'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.10%  '
$ws.Range('E24').Value = '  +0.07%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.807'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -4.74%  '
$ws.Range('B26').Value = 'Stellar'
$ws.Range('C26').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1284'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +4.07%  '
$ws.Range('B27').Value = 'Monero'
$ws.Range('C27').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '141.27'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.73%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.740'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.08%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.43'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.31%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.235'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.69%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.04873'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.11%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.303'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.30%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.185'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.84%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.558'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.80%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.372'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.28%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.8943'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.96%  '
$ws.Range('D37').Value = '1.128.23'
$ws.Range('E37').Value = '  +0.99%  '
$ws.Range('B38').Value = 'ImmutableX'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5497'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.77%  '
$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.532'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.41%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01560'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.07%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9976'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.33%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.591'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.20%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.7943'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.81%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '97.24'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.49%  '
$ws.Range('D45').Value = '1.771.30'
$ws.Range('E45').Value = '  -0.61%  '
$ws.Range('E46').Value = '  -7.89%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.4421'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.52%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '54.80'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.02%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05067'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.96%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.589'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.86%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.005'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.31%  '
